# Updated TPM-derived NATMI edge statistics (Hras -> Agtr1a) for rows 2-19,
# columns G:T: Ligand/Receptor average & total expression, derived
# specificities, and edge weights/specificities recomputed from new TPM input.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 13.93060933333333
$ws.Cells.Item(2, 8).Value = 41.791828
$ws.Cells.Item(2, 9).Value = 0.2372898381934647
$ws.Cells.Item(2, 10).Value = 0.2372898381934647
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 5.922617666666667
$ws.Cells.Item(2, 14).Value = 17.767853
$ws.Cells.Item(2, 15).Value = 0.2173916203328182
$ws.Cells.Item(2, 16).Value = 0.2173916203328182
$ws.Cells.Item(2, 17).Value = 82.50567294503153
$ws.Cells.Item(2, 18).Value = 742.5510565052839
$ws.Cells.Item(2, 19).Value = 0.05158482241338954
$ws.Cells.Item(2, 20).Value = 0.05158482241338953
# Row 3
$ws.Cells.Item(3, 7).Value = 13.93060933333333
$ws.Cells.Item(3, 8).Value = 41.791828
$ws.Cells.Item(3, 9).Value = 0.2372898381934647
$ws.Cells.Item(3, 10).Value = 0.2372898381934647
$ws.Cells.Item(3, 13).Value = 9.221608999999999
$ws.Cells.Item(3, 15).Value = 0.3384821772083041
$ws.Cells.Item(3, 16).Value = 0.3384821772083041
$ws.Cells.Item(3, 17).Value = 128.4626324037506
$ws.Cells.Item(3, 18).Value = 1156.163691633756
$ws.Cells.Item(3, 19).Value = 0.08031838106113011
$ws.Cells.Item(3, 20).Value = 0.08031838106113012
# Row 4
$ws.Cells.Item(4, 7).Value = 13.93060933333333
$ws.Cells.Item(4, 8).Value = 41.791828
$ws.Cells.Item(4, 9).Value = 0.2372898381934647
$ws.Cells.Item(4, 10).Value = 0.2372898381934647
$ws.Cells.Item(4, 13).Value = 12.099775
$ws.Cells.Item(4, 14).Value = 36.299325
$ws.Cells.Item(4, 15).Value = 0.4441262024588777
$ws.Cells.Item(4, 16).Value = 0.4441262024588777
$ws.Cells.Item(4, 17).Value = 168.5572385462333
$ws.Cells.Item(4, 18).Value = 1517.0151469161
$ws.Cells.Item(4, 19).Value = 0.105386634718945
$ws.Cells.Item(4, 20).Value = 0.105386634718945
# Row 5
$ws.Cells.Item(5, 9).Value = 0.1624178407807704
$ws.Cells.Item(5, 10).Value = 0.1624178407807703
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 5.922617666666667
$ws.Cells.Item(5, 14).Value = 17.767853
$ws.Cells.Item(5, 15).Value = 0.2173916203328182
$ws.Cells.Item(5, 16).Value = 0.2173916203328182
$ws.Cells.Item(5, 17).Value = 56.47268064202133
$ws.Cells.Item(5, 18).Value = 508.2541257781919
$ws.Cells.Item(5, 19).Value = 0.03530827757828935
$ws.Cells.Item(5, 20).Value = 0.03530827757828935
# Row 6
$ws.Cells.Item(6, 9).Value = 0.1624178407807704
$ws.Cells.Item(6, 10).Value = 0.1624178407807703
$ws.Cells.Item(6, 13).Value = 9.221608999999999
$ws.Cells.Item(6, 15).Value = 0.3384821772083041
$ws.Cells.Item(6, 16).Value = 0.3384821772083041
$ws.Cells.Item(6, 17).Value = 87.928853316592
$ws.Cells.Item(6, 18).Value = 791.3596798493279
$ws.Cells.Item(6, 19).Value = 0.05497554436494684
$ws.Cells.Item(6, 20).Value = 0.05497554436494684
# Row 7
$ws.Cells.Item(7, 9).Value = 0.1624178407807704
$ws.Cells.Item(7, 10).Value = 0.1624178407807703
$ws.Cells.Item(7, 13).Value = 12.099775
$ws.Cells.Item(7, 14).Value = 36.299325
$ws.Cells.Item(7, 15).Value = 0.4441262024588777
$ws.Cells.Item(7, 16).Value = 0.4441262024588777
$ws.Cells.Item(7, 17).Value = 115.3724194052
$ws.Cells.Item(7, 18).Value = 1038.3517746468
$ws.Cells.Item(7, 19).Value = 0.07213401883753419
$ws.Cells.Item(7, 20).Value = 0.07213401883753418
# Row 8
$ws.Cells.Item(8, 7).Value = 10.84369266666667
$ws.Cells.Item(8, 8).Value = 32.531078
$ws.Cells.Item(8, 9).Value = 0.1847082217815162
$ws.Cells.Item(8, 10).Value = 0.1847082217815162
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 5.922617666666667
$ws.Cells.Item(8, 14).Value = 17.767853
$ws.Cells.Item(8, 15).Value = 0.2173916203328182
$ws.Cells.Item(8, 16).Value = 0.2173916203328182
$ws.Cells.Item(8, 17).Value = 64.22304575950378
$ws.Cells.Item(8, 18).Value = 578.007411835534
$ws.Cells.Item(8, 19).Value = 0.04015401962187736
$ws.Cells.Item(8, 20).Value = 0.04015401962187735
# Row 9
$ws.Cells.Item(9, 7).Value = 10.84369266666667
$ws.Cells.Item(9, 8).Value = 32.531078
$ws.Cells.Item(9, 9).Value = 0.1847082217815162
$ws.Cells.Item(9, 10).Value = 0.1847082217815162
$ws.Cells.Item(9, 13).Value = 9.221608999999999
$ws.Cells.Item(9, 15).Value = 0.3384821772083041
$ws.Cells.Item(9, 16).Value = 0.3384821772083041
$ws.Cells.Item(9, 17).Value = 99.99629388816733
$ws.Cells.Item(9, 18).Value = 899.966644993506
$ws.Cells.Item(9, 19).Value = 0.06252044105688191
$ws.Cells.Item(9, 20).Value = 0.06252044105688191
# Row 10
$ws.Cells.Item(10, 7).Value = 10.84369266666667
$ws.Cells.Item(10, 8).Value = 32.531078
$ws.Cells.Item(10, 9).Value = 0.1847082217815162
$ws.Cells.Item(10, 10).Value = 0.1847082217815162
$ws.Cells.Item(10, 13).Value = 12.099775
$ws.Cells.Item(10, 14).Value = 36.299325
$ws.Cells.Item(10, 15).Value = 0.4441262024588777
$ws.Cells.Item(10, 16).Value = 0.4441262024588777
$ws.Cells.Item(10, 17).Value = 131.2062414358167
$ws.Cells.Item(10, 18).Value = 1180.85617292235
$ws.Cells.Item(10, 19).Value = 0.08203376110275695
$ws.Cells.Item(10, 20).Value = 0.08203376110275695
# Row 11
$ws.Cells.Item(11, 7).Value = 8.514172333333333
$ws.Cells.Item(11, 8).Value = 25.542517
$ws.Cells.Item(11, 9).Value = 0.1450278682708931
$ws.Cells.Item(11, 10).Value = 0.1450278682708931
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 5.922617666666667
$ws.Cells.Item(11, 14).Value = 17.767853
$ws.Cells.Item(11, 15).Value = 0.2173916203328182
$ws.Cells.Item(11, 16).Value = 0.2173916203328182
$ws.Cells.Item(11, 17).Value = 50.42618747844455
$ws.Cells.Item(11, 18).Value = 453.8356873060009
$ws.Cells.Item(11, 19).Value = 0.03152784327682396
$ws.Cells.Item(11, 20).Value = 0.03152784327682395
# Row 12
$ws.Cells.Item(12, 7).Value = 8.514172333333333
$ws.Cells.Item(12, 8).Value = 25.542517
$ws.Cells.Item(12, 9).Value = 0.1450278682708931
$ws.Cells.Item(12, 10).Value = 0.1450278682708931
$ws.Cells.Item(12, 13).Value = 9.221608999999999
$ws.Cells.Item(12, 15).Value = 0.3384821772083041
$ws.Cells.Item(12, 16).Value = 0.3384821772083041
$ws.Cells.Item(12, 17).Value = 78.51436821661765
$ws.Cells.Item(12, 18).Value = 706.6293139495589
$ws.Cells.Item(12, 19).Value = 0.04908934860821101
$ws.Cells.Item(12, 20).Value = 0.04908934860821101
# Row 13
$ws.Cells.Item(13, 7).Value = 8.514172333333333
$ws.Cells.Item(13, 8).Value = 25.542517
$ws.Cells.Item(13, 9).Value = 0.1450278682708931
$ws.Cells.Item(13, 10).Value = 0.1450278682708931
$ws.Cells.Item(13, 13).Value = 12.099775
$ws.Cells.Item(13, 14).Value = 36.299325
$ws.Cells.Item(13, 15).Value = 0.4441262024588777
$ws.Cells.Item(13, 16).Value = 0.4441262024588777
$ws.Cells.Item(13, 17).Value = 103.0195695445583
$ws.Cells.Item(13, 18).Value = 927.176125901025
$ws.Cells.Item(13, 19).Value = 0.0644106763858581
$ws.Cells.Item(13, 20).Value = 0.0644106763858581
# Row 14
$ws.Cells.Item(14, 7).Value = 6.413260666666666
$ws.Cells.Item(14, 8).Value = 19.239782
$ws.Cells.Item(14, 9).Value = 0.1092415665009325
$ws.Cells.Item(14, 10).Value = 0.1092415665009325
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 5.922617666666667
$ws.Cells.Item(14, 14).Value = 17.767853
$ws.Cells.Item(14, 15).Value = 0.2173916203328182
$ws.Cells.Item(14, 16).Value = 0.2173916203328182
$ws.Cells.Item(14, 17).Value = 37.98329092533844
$ws.Cells.Item(14, 18).Value = 341.8496183280459
$ws.Cells.Item(14, 19).Value = 0.02374820114933304
$ws.Cells.Item(14, 20).Value = 0.02374820114933304
# Row 15
$ws.Cells.Item(15, 7).Value = 6.413260666666666
$ws.Cells.Item(15, 8).Value = 19.239782
$ws.Cells.Item(15, 9).Value = 0.1092415665009325
$ws.Cells.Item(15, 10).Value = 0.1092415665009325
$ws.Cells.Item(15, 13).Value = 9.221608999999999
$ws.Cells.Item(15, 15).Value = 0.3384821772083041
$ws.Cells.Item(15, 16).Value = 0.3384821772083041
$ws.Cells.Item(15, 17).Value = 59.14058228307932
$ws.Cells.Item(15, 18).Value = 532.2652405477139
$ws.Cells.Item(15, 19).Value = 0.03697632327088138
$ws.Cells.Item(15, 20).Value = 0.03697632327088139
# Row 16
$ws.Cells.Item(16, 7).Value = 6.413260666666666
$ws.Cells.Item(16, 8).Value = 19.239782
$ws.Cells.Item(16, 9).Value = 0.1092415665009325
$ws.Cells.Item(16, 10).Value = 0.1092415665009325
$ws.Cells.Item(16, 13).Value = 12.099775
$ws.Cells.Item(16, 14).Value = 36.299325
$ws.Cells.Item(16, 15).Value = 0.4441262024588777
$ws.Cells.Item(16, 16).Value = 0.4441262024588777
$ws.Cells.Item(16, 17).Value = 77.59901108301666
$ws.Cells.Item(16, 18).Value = 698.39109974715
$ws.Cells.Item(16, 19).Value = 0.04851704208071811
$ws.Cells.Item(16, 20).Value = 0.04851704208071811
# Row 17
$ws.Cells.Item(17, 7).Value = 9.470323666666667
$ws.Cells.Item(17, 8).Value = 28.410971
$ws.Cells.Item(17, 9).Value = 0.1613146644724231
$ws.Cells.Item(17, 10).Value = 0.1613146644724231
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 5.922617666666667
$ws.Cells.Item(17, 14).Value = 17.767853
$ws.Cells.Item(17, 15).Value = 0.2173916203328182
$ws.Cells.Item(17, 16).Value = 0.2173916203328182
$ws.Cells.Item(17, 17).Value = 56.08910625725144
$ws.Cells.Item(17, 18).Value = 504.8019563152631
$ws.Cells.Item(17, 19).Value = 0.03506845629310497
$ws.Cells.Item(17, 20).Value = 0.03506845629310496
# Row 18
$ws.Cells.Item(18, 7).Value = 9.470323666666667
$ws.Cells.Item(18, 8).Value = 28.410971
$ws.Cells.Item(18, 9).Value = 0.1613146644724231
$ws.Cells.Item(18, 10).Value = 0.1613146644724231
$ws.Cells.Item(18, 13).Value = 9.221608999999999
$ws.Cells.Item(18, 15).Value = 0.3384821772083041
$ws.Cells.Item(18, 16).Value = 0.3384821772083041
$ws.Cells.Item(18, 17).Value = 87.33162195744633
$ws.Cells.Item(18, 18).Value = 785.984597617017
$ws.Cells.Item(18, 19).Value = 0.05460213884625284
$ws.Cells.Item(18, 20).Value = 0.05460213884625285
# Row 19
$ws.Cells.Item(19, 7).Value = 9.470323666666667
$ws.Cells.Item(19, 8).Value = 28.410971
$ws.Cells.Item(19, 9).Value = 0.1613146644724231
$ws.Cells.Item(19, 10).Value = 0.1613146644724231
$ws.Cells.Item(19, 13).Value = 12.099775
$ws.Cells.Item(19, 14).Value = 36.299325
$ws.Cells.Item(19, 15).Value = 0.4441262024588777
$ws.Cells.Item(19, 16).Value = 0.4441262024588777
$ws.Cells.Item(19, 17).Value = 114.5887855438417
$ws.Cells.Item(19, 18).Value = 1031.299069894575
$ws.Cells.Item(19, 19).Value = 0.07164406933306532
$ws.Cells.Item(19, 20).Value = 0.07164406933306532
